$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Good Morning" greeting text in E8 to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the active cell selection to E8
$ws.Range("E8").Select()
